# "update flow voucher green" - add two new API-master rows:
#   row 4: booking/checkin  -> CheckInBagPlayer1VoucherTest
#   row 5: booking/edit_booking -> EditBooking1PlayerVoucherTest

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 4 : booking/checkin ----
$ws.Range("A4").Value = "booking/checkin"
$ws.Range("B4").Value = "CheckInBagPlayer1VoucherTest"
$ws.Range("C4").Value = "checkInBagData"
$ws.Range("D4").Value = "BASE_URL"
$ws.Range("E4").Value = "/golf-cms/api/booking/check-in"
$ws.Range("F4").Value = "POST"
$ws.Range("G4").Value = "JSON"
$ws.Range("H4").Value = "input_excel_file/booking/CheckIn.xlsx"
$ws.Range("I4").Value = "Check_In_Bag_Player1_VC"
$ws.Range("J4").Value = "input_json_file/booking/check_in/"
$ws.Range("K4").Value = "testCheckInBag"
$ws.Range("L4").Value = "check_in_bag_player1_voucher_case_id"

# ---- Row 5 : booking/edit_booking ----
$ws.Range("A5").Value = "booking/edit_booking"
$ws.Range("B5").Value = "EditBooking1PlayerVoucherTest"
$ws.Range("C5").Value = "EditBookingData"
$ws.Range("D5").Value = "BASE_URL"
$ws.Range("E5").Value = "/golf-cms/api/booking/update"
$ws.Range("F5").Value = "POST"
$ws.Range("G5").Value = "JSON"
$ws.Range("H5").Value = "input_excel_file/booking/Create_Booking_Batch.xlsx"
$ws.Range("I5").Value = "Edit_Booking_1_Player_VC"
$ws.Range("J5").Value = "input_json_file/booking/edit_booking/"
$ws.Range("K5").Value = "testEditBooking"
$ws.Range("L5").Value = "edit_booking_1player_voucher_case_id"

# New rows repeat the wrapped / vertically centered look of row 3 and
# need the same fixed 45pt row height.
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(5).RowHeight = 45

# Column widths were retuned for the new content.
$ws.Columns.Item(1).ColumnWidth = 22.0
$ws.Columns.Item(4).ColumnWidth = 12.666666666666666
$ws.Columns.Item(8).ColumnWidth = 24.333333333333336
$ws.Columns.Item(11).ColumnWidth = 20.0
$ws.Columns.Item(12).ColumnWidth = 22.333333333333336

# The workbook's outline-level watermark moves from 2 to 4 (rows were
# grouped up to that depth at some point); bump it without leaving any
# row with a visible outline level by grouping/ungrouping a scratch row.
$ws.Rows.Item(1000).OutlineLevel = 4
$ws.Rows.Item(1000).Delete()

# Selection moves to E10.
$null = $ws.Range("E10").Select()
